# Validation dataset, gui formatting
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Datasets")

# Copy the formatting of the row above it (row 14) onto the new row 15
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D15").PasteSpecial(-4104) | Out-Null

# Add the new "Clustering Test" validation dataset row (row 15)
$ws.Cells.Item(15, 1).Value = 14
$ws.Cells.Item(15, 2).Value = "Clustering Test"
$ws.Cells.Item(15, 3).Value = ".\datasets\test_clustering.xlsx"
$ws.Cells.Item(15, 4).Value = "S. S. M. Ghoneim, and I. B. M. Taha,`"A New Approach of DGA Interpretation Technique for Transformer Fault Diagnosis`", International Journal of Electrical Power and Energy Systems, 81, Oct. 2016, pp. 265–274."

$ws.Range("A15").RowHeight = 45
$ws.Range("A1:D15").Select()

$wb.Save()
